$d = $word.ActiveDocument

# The first paragraph currently holds two runs: "**ID__AFFARS_mp_5332_7_topic_2__ID**"
# followed by a lone space run. Replace the whole paragraph's text (excluding its
# end-of-paragraph mark) with the new single-run ID placeholder, which both renames
# the ID and collapses the two runs into one.
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$r.End = $r.End - 1
$r.Text = "**ID__AFFARS_RELEASE_OF_SOLICITATIONS__ID**"

# Give the paragraph a thin paragraph border (5 twips on every side) and widen its
# left indent from 120 to 225 twips (= 11.25pt, since Word's indent properties are
# expressed in points while ind/@w:left is stored in twentieths of a point).
$p1.Format.LeftIndent = 225 / 20.0
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5
